$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data to the log
$ws.Range("A23").Value = 45643
$ws.Range("A23").NumberFormat = "d-mmm"
$ws.Range("B23").Value = "puzzle design and lore work"
$ws.Range("C23").Value = 5

# Update the selected cell to match the new state
$ws.Range("E21").Select()
